$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 355, shifting existing rows 355:440 down to 356:441
$ws.Rows(355).Insert()

# Populate the newly inserted row 355 with the new weekly data point
$ws.Cells.Item(355, 1).Value = 3
$ws.Cells.Item(355, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(355, 3).Value = "Coquimbo"
$ws.Cells.Item(355, 4).Value = 44785
$ws.Cells.Item(355, 5).Value = 5
$ws.Cells.Item(355, 6).Value = 100112017
$ws.Cells.Item(355, 7).Value = "Apio"
$ws.Cells.Item(355, 8).Value = "Americana (o)"
$ws.Cells.Item(355, 9).Value = "Primera"
$ws.Cells.Item(355, 10).Value = 230
$ws.Cells.Item(355, 11).Value = 9000
$ws.Cells.Item(355, 12).Value = 9500
$ws.Cells.Item(355, 13).Value = 9239
$ws.Cells.Item(355, 14).Value = "$/docena de matas"
$ws.Cells.Item(355, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(355, 16).Value = 1540
$ws.Cells.Item(355, 17).Value = 6
$ws.Cells.Item(355, 18).Value = "Hortaliza"
